$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list (Price / Volume(1h)) - GitHub Actions refresh
# Column D (Price) values that look like plain numbers need a leading apostrophe
# so Excel keeps them as text (matching the source feed's formatted strings)
# instead of silently parsing them into floating point numbers.

$ws.Cells.Item(2, 4).Value = "26.179.83"
$ws.Cells.Item(2, 5).Value = "  -0.61%  "
$ws.Cells.Item(3, 4).Value = "1.586.20"
$ws.Cells.Item(3, 5).Value = "  -0.40%  "
$ws.Cells.Item(4, 5).Value = "  +0.02%  "
$ws.Cells.Item(5, 4).Value = "'211.87"
$ws.Cells.Item(5, 5).Value = "  +0.78%  "
$ws.Cells.Item(6, 5).Value = "  -0.41%  "
$ws.Cells.Item(7, 5).Value = "  +0.02%  "
$ws.Cells.Item(8, 5).Value = "  -0.38%  "
$ws.Cells.Item(9, 4).Value = "'0.0603"
$ws.Cells.Item(9, 5).Value = "  -1.39%  "
$ws.Cells.Item(10, 4).Value = "'19.19"
$ws.Cells.Item(10, 5).Value = "  -2.10%  "
$ws.Cells.Item(11, 4).Value = "'0.0846"
$ws.Cells.Item(11, 5).Value = "  +0.19%  "
$ws.Cells.Item(12, 4).Value = "1.811.13"
$ws.Cells.Item(12, 5).Value = "  -0.26%  "
$ws.Cells.Item(13, 4).Value = "1.598.57"
$ws.Cells.Item(13, 5).Value = "  -0.35%  "
$ws.Cells.Item(14, 5).Value = "  -1.66%  "
$ws.Cells.Item(15, 4).Value = "'0.516"
$ws.Cells.Item(15, 5).Value = "  -0.45%  "
$ws.Cells.Item(16, 4).Value = "'63.90"
$ws.Cells.Item(16, 5).Value = "  -1.22%  "
$ws.Cells.Item(17, 4).Value = "26.186.87"
$ws.Cells.Item(17, 5).Value = "  -0.53%  "
$ws.Cells.Item(18, 5).Value = "  -0.82%  "
$ws.Cells.Item(19, 4).Value = "'214.00"
$ws.Cells.Item(19, 5).Value = "  +0.97%  "
$ws.Cells.Item(20, 5).Value = "  -2.86%  "
$ws.Cells.Item(22, 5).Value = "  -0.94%  "
$ws.Cells.Item(23, 4).Value = "'2.17"
$ws.Cells.Item(23, 5).Value = "  -1.12%  "
$ws.Cells.Item(24, 4).Value = "'8.95"
$ws.Cells.Item(24, 5).Value = "  +0.31%  "
$ws.Cells.Item(25, 4).Value = "'143.89"
$ws.Cells.Item(25, 5).Value = "  -1.05%  "
$ws.Cells.Item(26, 5).Value = "  +0.04%  "
$ws.Cells.Item(27, 4).Value = "'6.99"
$ws.Cells.Item(27, 5).Value = "  -0.86%  "
$ws.Cells.Item(28, 5).Value = "  -1.07%  "
$ws.Cells.Item(29, 4).Value = "'15.07"
$ws.Cells.Item(29, 5).Value = "  -1.34%  "
$ws.Cells.Item(30, 5).Value = "  -1.73%  "
$ws.Cells.Item(31, 5).Value = "  +0.47%  "
$ws.Cells.Item(32, 4).Value = "'3.18"
$ws.Cells.Item(32, 5).Value = "  -1.02%  "
$ws.Cells.Item(33, 4).Value = "1.391.80"
$ws.Cells.Item(33, 5).Value = "  +6.88%  "
$ws.Cells.Item(34, 4).Value = "'2.93"
$ws.Cells.Item(34, 5).Value = "  -1.81%  "
$ws.Cells.Item(35, 5).Value = "  -0.52%  "
$ws.Cells.Item(36, 5).Value = "  -1.46%  "
$ws.Cells.Item(37, 4).Value = "'0.585"
$ws.Cells.Item(37, 5).Value = "  -4.72%  "
$ws.Cells.Item(38, 5).Value = "  -0.96%  "
$ws.Cells.Item(39, 5).Value = "  +0.54%  "
$ws.Cells.Item(40, 5).Value = "  +3.99%  "
$ws.Cells.Item(41, 5).Value = "  -0.05%  "
$ws.Cells.Item(42, 4).Value = "'0.939"
$ws.Cells.Item(42, 5).Value = "  -14.79%  "
$ws.Cells.Item(43, 4).Value = "'0.766"
$ws.Cells.Item(43, 5).Value = "  +0.60%  "
$ws.Cells.Item(44, 5).Value = "  +0.03%  "
$ws.Cells.Item(45, 4).Value = "1.722.87"
$ws.Cells.Item(45, 5).Value = "  -0.30%  "
$ws.Cells.Item(46, 4).Value = "'60.94"
$ws.Cells.Item(46, 5).Value = "  -2.81%  "
$ws.Cells.Item(47, 4).Value = "'85.95"
$ws.Cells.Item(47, 5).Value = "  -2.91%  "
$ws.Cells.Item(48, 4).Value = "'1.48"
$ws.Cells.Item(48, 5).Value = "  -1.83%  "
$ws.Cells.Item(49, 4).Value = "'0.0974"
$ws.Cells.Item(49, 5).Value = "  -1.16%  "
$ws.Cells.Item(50, 4).Value = "'0.0499"
$ws.Cells.Item(50, 5).Value = "  -1.19%  "
$ws.Cells.Item(51, 4).Value = "'0.999"
$ws.Cells.Item(51, 5).Value = "  -0.08%  "
